$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 'FAPs'
$ws.Cells.Item(2, 7).Value = 61.526261
$ws.Cells.Item(2, 8).Value = 184.578783
$ws.Cells.Item(2, 9).Value = 0.9684015388399598
$ws.Cells.Item(2, 10).Value = 0.9684015388399598
$ws.Cells.Item(2, 11).Value = 3.0
$ws.Cells.Item(2, 12).Value = 1.0
$ws.Cells.Item(2, 13).Value = 2.214957333333333
$ws.Cells.Item(2, 14).Value = 6.644871999999999
$ws.Cells.Item(2, 15).Value = 0.8812411509483107
$ws.Cells.Item(2, 16).Value = 0.8812411509483107
$ws.Cells.Item(2, 17).Value = 136.2780429945306
$ws.Cells.Item(2, 18).Value = 1226.502386950776
$ws.Cells.Item(2, 19).Value = 0.8533952866674414
$ws.Cells.Item(2, 20).Value = 0.8533952866674414
$ws.Cells.Item(3, 4).Value = 'MuSCs'
$ws.Cells.Item(3, 7).Value = 61.526261
$ws.Cells.Item(3, 8).Value = 184.578783
$ws.Cells.Item(3, 9).Value = 0.9684015388399598
$ws.Cells.Item(3, 10).Value = 0.9684015388399598
$ws.Cells.Item(3, 13).Value = 0.274148
$ws.Cells.Item(3, 14).Value = 0.8224440000000001
$ws.Cells.Item(3, 15).Value = 0.1090723037479928
$ws.Cells.Item(3, 16).Value = 0.1090723037479928
$ws.Cells.Item(3, 17).Value = 16.867301400628
$ws.Cells.Item(3, 18).Value = 151.805712605652
$ws.Cells.Item(3, 19).Value = 0.1056257867943758
$ws.Cells.Item(3, 20).Value = 0.1056257867943758
$ws.Cells.Item(4, 4).Value = 'Resolving-Mac'
$ws.Cells.Item(4, 7).Value = 61.526261
$ws.Cells.Item(4, 8).Value = 184.578783
$ws.Cells.Item(4, 9).Value = 0.9684015388399598
$ws.Cells.Item(4, 10).Value = 0.9684015388399598
$ws.Cells.Item(4, 11).Value = 1.0
$ws.Cells.Item(4, 12).Value = 0.3333333333333333
$ws.Cells.Item(4, 13).Value = 0.02434666666666667
$ws.Cells.Item(4, 14).Value = 0.07304
$ws.Cells.Item(4, 15).Value = 0.009686545303696538
$ws.Cells.Item(4, 16).Value = 0.009686545303696536
$ws.Cells.Item(4, 17).Value = 1.497959367813333
$ws.Cells.Item(4, 18).Value = 13.48163431032
$ws.Cells.Item(4, 19).Value = 0.009380465378142713
$ws.Cells.Item(4, 20).Value = 0.009380465378142711
$ws.Cells.Item(5, 4).Value = 'FAPs'
$ws.Cells.Item(5, 9).Value = 0.002979850677668077
$ws.Cells.Item(5, 10).Value = 0.002979850677668078
$ws.Cells.Item(5, 11).Value = 3.0
$ws.Cells.Item(5, 12).Value = 1.0
$ws.Cells.Item(5, 13).Value = 2.214957333333333
$ws.Cells.Item(5, 14).Value = 6.644871999999999
$ws.Cells.Item(5, 15).Value = 0.8812411509483107
$ws.Cells.Item(5, 16).Value = 0.8812411509483107
$ws.Cells.Item(5, 17).Value = 0.4193386756231111
$ws.Cells.Item(5, 18).Value = 3.774048080608
$ws.Cells.Item(5, 19).Value = 0.00262596704084232
$ws.Cells.Item(5, 20).Value = 0.00262596704084232
$ws.Cells.Item(6, 4).Value = 'MuSCs'
$ws.Cells.Item(6, 9).Value = 0.002979850677668077
$ws.Cells.Item(6, 10).Value = 0.002979850677668078
$ws.Cells.Item(6, 13).Value = 0.274148
$ws.Cells.Item(6, 14).Value = 0.8224440000000001
$ws.Cells.Item(6, 15).Value = 0.1090723037479928
$ws.Cells.Item(6, 16).Value = 0.1090723037479928
$ws.Cells.Item(6, 17).Value = 0.05190206489066667
$ws.Cells.Item(6, 18).Value = 0.467118584016
$ws.Cells.Item(6, 19).Value = 0.0003250191782382747
$ws.Cells.Item(6, 20).Value = 0.0003250191782382748
$ws.Cells.Item(7, 4).Value = 'Resolving-Mac'
$ws.Cells.Item(7, 9).Value = 0.002979850677668077
$ws.Cells.Item(7, 10).Value = 0.002979850677668078
$ws.Cells.Item(7, 11).Value = 1.0
$ws.Cells.Item(7, 12).Value = 0.3333333333333333
$ws.Cells.Item(7, 13).Value = 0.02434666666666667
$ws.Cells.Item(7, 14).Value = 0.07304
$ws.Cells.Item(7, 15).Value = 0.009686545303696538
$ws.Cells.Item(7, 16).Value = 0.009686545303696536
$ws.Cells.Item(7, 17).Value = 0.004609343395555555
$ws.Cells.Item(7, 18).Value = 0.04148409056
$ws.Cells.Item(7, 19).Value = 0.00002886445858748266
$ws.Cells.Item(7, 20).Value = 0.00002886445858748266
$ws.Cells.Item(8, 4).Value = 'FAPs'
$ws.Cells.Item(8, 7).Value = 0.6472316666666668
$ws.Cells.Item(8, 8).Value = 1.941695
$ws.Cells.Item(8, 9).Value = 0.01018719700821657
$ws.Cells.Item(8, 10).Value = 0.01018719700821657
$ws.Cells.Item(8, 11).Value = 3.0
$ws.Cells.Item(8, 12).Value = 1.0
$ws.Cells.Item(8, 13).Value = 2.214957333333333
$ws.Cells.Item(8, 14).Value = 6.644871999999999
$ws.Cells.Item(8, 15).Value = 0.8812411509483107
$ws.Cells.Item(8, 16).Value = 0.8812411509483107
$ws.Cells.Item(8, 17).Value = 1.433590526448889
$ws.Cells.Item(8, 18).Value = 12.90231473804
$ws.Cells.Item(8, 19).Value = 0.00897737721645796
$ws.Cells.Item(8, 20).Value = 0.00897737721645796
$ws.Cells.Item(9, 4).Value = 'MuSCs'
$ws.Cells.Item(9, 7).Value = 0.6472316666666668
$ws.Cells.Item(9, 8).Value = 1.941695
$ws.Cells.Item(9, 9).Value = 0.01018719700821657
$ws.Cells.Item(9, 10).Value = 0.01018719700821657
$ws.Cells.Item(9, 13).Value = 0.274148
$ws.Cells.Item(9, 14).Value = 0.8224440000000001
$ws.Cells.Item(9, 15).Value = 0.1090723037479928
$ws.Cells.Item(9, 16).Value = 0.1090723037479928
$ws.Cells.Item(9, 17).Value = 0.1774372669533334
$ws.Cells.Item(9, 18).Value = 1.59693540258
$ws.Cells.Item(9, 19).Value = 0.001111141046420842
$ws.Cells.Item(9, 20).Value = 0.001111141046420842
$ws.Cells.Item(10, 4).Value = 'Resolving-Mac'
$ws.Cells.Item(10, 7).Value = 0.6472316666666668
$ws.Cells.Item(10, 8).Value = 1.941695
$ws.Cells.Item(10, 9).Value = 0.01018719700821657
$ws.Cells.Item(10, 10).Value = 0.01018719700821657
$ws.Cells.Item(10, 11).Value = 1.0
$ws.Cells.Item(10, 12).Value = 0.3333333333333333
$ws.Cells.Item(10, 13).Value = 0.02434666666666667
$ws.Cells.Item(10, 14).Value = 0.07304
$ws.Cells.Item(10, 15).Value = 0.009686545303696538
$ws.Cells.Item(10, 16).Value = 0.009686545303696536
$ws.Cells.Item(10, 17).Value = 0.01575793364444445
$ws.Cells.Item(10, 18).Value = 0.1418214028
$ws.Cells.Item(10, 19).Value = 0.00009867874533777168
$ws.Cells.Item(10, 20).Value = 0.00009867874533777165
$ws.Cells.Item(11, 4).Value = 'FAPs'
$ws.Cells.Item(11, 7).Value = 1.171018333333333
$ws.Cells.Item(11, 8).Value = 3.513055
$ws.Cells.Item(11, 9).Value = 0.01843141347415545
$ws.Cells.Item(11, 10).Value = 0.01843141347415545
$ws.Cells.Item(11, 11).Value = 3.0
$ws.Cells.Item(11, 12).Value = 1.0
$ws.Cells.Item(11, 13).Value = 2.214957333333333
$ws.Cells.Item(11, 14).Value = 6.644871999999999
$ws.Cells.Item(11, 15).Value = 0.8812411509483107
$ws.Cells.Item(11, 16).Value = 0.8812411509483107
$ws.Cells.Item(11, 17).Value = 2.593755644884444
$ws.Cells.Item(11, 18).Value = 23.34380080396
$ws.Cells.Item(11, 19).Value = 0.01624252002356895
$ws.Cells.Item(11, 20).Value = 0.01624252002356895
$ws.Cells.Item(12, 4).Value = 'MuSCs'
$ws.Cells.Item(12, 7).Value = 1.171018333333333
$ws.Cells.Item(12, 8).Value = 3.513055
$ws.Cells.Item(12, 9).Value = 0.01843141347415545
$ws.Cells.Item(12, 10).Value = 0.01843141347415545
$ws.Cells.Item(12, 13).Value = 0.274148
$ws.Cells.Item(12, 14).Value = 0.8224440000000001
$ws.Cells.Item(12, 15).Value = 0.1090723037479928
$ws.Cells.Item(12, 16).Value = 0.1090723037479928
$ws.Cells.Item(12, 17).Value = 0.3210323340466666
$ws.Cells.Item(12, 18).Value = 2.88929100642
$ws.Cells.Item(12, 19).Value = 0.002010356728957931
$ws.Cells.Item(12, 20).Value = 0.002010356728957931
$ws.Cells.Item(13, 4).Value = 'Resolving-Mac'
$ws.Cells.Item(13, 7).Value = 1.171018333333333
$ws.Cells.Item(13, 8).Value = 3.513055
$ws.Cells.Item(13, 9).Value = 0.01843141347415545
$ws.Cells.Item(13, 10).Value = 0.01843141347415545
$ws.Cells.Item(13, 11).Value = 1.0
$ws.Cells.Item(13, 12).Value = 0.3333333333333333
$ws.Cells.Item(13, 13).Value = 0.02434666666666667
$ws.Cells.Item(13, 14).Value = 0.07304
$ws.Cells.Item(13, 15).Value = 0.009686545303696538
$ws.Cells.Item(13, 16).Value = 0.009686545303696536
$ws.Cells.Item(13, 17).Value = 0.02851039302222222
$ws.Cells.Item(13, 18).Value = 0.2565935372
$ws.Cells.Item(13, 19).Value = 0.0001785367216285696
$ws.Cells.Item(13, 20).Value = 0.0001785367216285696

Write-Host "Applied changes"
